$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("List_1")
$ws2 = $wb.Worksheets.Item("List_2")

# ---------------------------------------------------------------------------
# 1) List_2 (sheet2): append a new data row (row 3) under the existing table,
#    then a trailing blank-but-formatted row (row 4) - the validator output.
# ---------------------------------------------------------------------------

# Bring formatting down from row 2 onto row 3 (reuses existing style indices
# instead of minting new ones).
[void]$ws2.Range("A2:Q2").Copy()
[void]$ws2.Range("A3:Q3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("A3").Value = "wilbert"
$ws2.Range("B3").Value = "jnroasvn"
$ws2.Range("C3").Value = "mrv"
$ws2.Range("D3").Value = 33216
$ws2.Range("E3").Value = 2
$ws2.Range("F3").Value = "WIRI"
$ws2.Range("G3").Value = "K"
$ws2.Range("H3").Value = "nutley"
$ws2.Range("I3").Value = 2
$ws2.Range("J3").Value = 2
$ws2.Range("K3").Value = 1
$ws2.Range("L3").Value = 6
$ws2.Range("M3").Value = 8
$ws2.Range("N3").Value = 8
$ws2.Range("O3").Value = "yes"
$ws2.Range("P3").Value = "ayes"
$ws2.Range("Q3").Value = "hmm"

# Trailing row 4: just a formatted (date-styled), empty D4 cell.
[void]$ws2.Range("D2").Copy()
[void]$ws2.Range("D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Grow the List_2 table to cover the new data row + the trailing blank row.
$lo2 = $ws2.ListObjects.Item(1)
[void]$lo2.Resize($ws2.Range("A1:Q4"))

[void]$ws2.Range("A4").Select()

# ---------------------------------------------------------------------------
# 2) List_1 (sheet1): append a new data row (row 3).
# ---------------------------------------------------------------------------

[void]$ws1.Range("A2:U2").Copy()
[void]$ws1.Range("A3:U3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A3").Value = "willis"
$ws1.Range("B3").Value = "dfnjfd"
$ws1.Range("C3").Value = "fjkndkj"
$ws1.Range("D3").Value = 36892
$ws1.Range("E3").Value = 2
$ws1.Range("F3").Value = "WIRT"
$ws1.Range("G3").Value = 38019
$ws1.Range("H3").Value = 39969
$ws1.Range("I3").Value = 40890
$ws1.Range("J3").Value = "K1"
$ws1.Range("K3").Value = "Place"
$ws1.Range("L3").Value = 2
$ws1.Range("M3").Value = "Place 2"
$ws1.Range("N3").Value = "Yes"
$ws1.Range("O3").Value = 20
$ws1.Range("P3").Value = 1
$ws1.Range("Q3").Value = 11
$ws1.Range("R3").Value = "yes"
$ws1.Range("S3").Value = "yes"
$ws1.Range("T3").Value = "yes"
$ws1.Range("U3").Value = "yes"

# List_1's table range already spans A1:U3 - nothing to resize there.

[void]$ws2.Activate()
[void]$ws2.Range("A4").Select()
[void]$ws1.Activate()
[void]$ws1.Range("U3").Select()
